$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Date placeholder fields (master, layouts, notes master) ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "2018. 1. 27.") {
        $sh.TextFrame.TextRange.Text = "2018. 3. 5."
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $lay = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "2018. 1. 27.") {
            $sh.TextFrame.TextRange.Text = "2018. 3. 5."
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "2018. 1. 27.") {
        $sh.TextFrame.TextRange.Text = "2018. 3. 5."
    }
}

# --- 2. Simple single-letter label textboxes: reposition/resize + rename ---
$sh = $s.Shapes.Item(12)
$sh.Left = 79.529766163713
$sh.Top = 61.086851347375
$sh.Width = 31.687796706102
$sh.Height = 29.081259842520
$sh.TextFrame.TextRange.Text = "lr"

$sh = $s.Shapes.Item(13)
$sh.Left = 95.258740157480
$sh.Top = 45.756063468963
$sh.Width = 34.920316867979
$sh.Height = 29.081259842520
$sh.TextFrame.TextRange.Text = "lc"

$sh = $s.Shapes.Item(14)
$sh.Left = 149.023307801870
$sh.Top = 45.756063468963
$sh.Width = 34.932600094127
$sh.Height = 29.081259842520
$sh.TextFrame.TextRange.Text = "rc"

$sh = $s.Shapes.Item(15)
$sh.Left = 135.485984967224
$sh.Top = 61.924017178543
$sh.Width = 32.934568598064
$sh.Height = 29.081259842520
$sh.TextFrame.TextRange.Text = "lc"

$sh = $s.Shapes.Item(17)
$sh.Left = 231.099769739992
$sh.Top = 62.225355284383
$sh.Width = 63.696222379790
$sh.Height = 31.504725601542
$sh.TextFrame.TextRange.Text = "lr"

$sh = $s.Shapes.Item(18)
$sh.Left = 230.402839890878
$sh.Top = 91.816692913386
$sh.Width = 63.791181579199
$sh.Height = 31.504725601542
$sh.TextFrame.TextRange.Text = "lc"

$sh = $s.Shapes.Item(19)
$sh.Left = 227.866692913386
$sh.Top = 118.556377952756
$sh.Width = 67.946929133858
$sh.Height = 31.504725601542
$sh.TextFrame.TextRange.Text = "rc"

# --- 3. Multi-run label textboxes (n*k*m style -> lr*lc*rc style) ---
$sh = $s.Shapes.Item(20)
$sh.Left = 217.129133858268
$sh.Top = 278.972677165354
$sh.TextFrame.TextRange.Text = "lr*lc*rc"
$sh.TextFrame.TextRange.Characters(1, 2).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 5).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 6).Font.Bold = $true

$sh = $s.Shapes.Item(21)
$sh.Left = 217.129133858268
$sh.Top = 214.623543307087
$sh.TextFrame.TextRange.Text = "lr*rc"
$sh.TextFrame.TextRange.Characters(1, 2).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true

$sh = $s.Shapes.Item(22)
$sh.Left = 218.236778329264
$sh.Top = 181.982677165354
$sh.TextFrame.TextRange.Text = "lc*rc"
$sh.TextFrame.TextRange.Characters(1, 2).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true

$sh = $s.Shapes.Item(23)
$sh.Left = 217.129055118110
$sh.Top = 149.736062992126
$sh.TextFrame.TextRange.Text = "lr*lc"
$sh.TextFrame.TextRange.Characters(1, 2).Font.Bold = $true
$sh.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true

# --- 4. Resize the roundRect frame that groups these labels ---
$sh = $s.Shapes.Item(85)
$sh.Left = 208.069133858268
$sh.Top = 38.117639702624
$sh.Width = 105.234333092847
$sh.Height = 301.554094488189

# --- 5. Add the new "lr*lc+lc*rc" textbox ---
$newBox = $s.Shapes.AddTextbox(1, 204.905906703904, 246.082440944882, 111.228977331627, 31.504725601542)
$newBox.TextFrame.TextRange.Text = "lr*lc+lc*rc"
$newBox.TextFrame.TextRange.Font.Size = 20
$newBox.TextFrame.TextRange.Font.Bold = $true
$newBox.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$newBox.TextFrame.TextRange.Characters(1, 2).Font.Bold = $true
$newBox.TextFrame.TextRange.Characters(1, 3).Font.Bold = $true
$newBox.TextFrame.TextRange.Characters(1, 8).Font.Bold = $true
$newBox.TextFrame.TextRange.Characters(1, 9).Font.Bold = $true

